# Auto-generated edit script applying the scheduled-runner price/profit updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 3799.4
$ws.Range("J2").Value = 6999.5
$ws.Range("L2").Value = 6999.5
$ws.Range("N2").Value = -7225.5

# Row 18
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

# Row 32
$ws.Range("H32").Value = 9164.833000000001
$ws.Range("J32").Value = 9497.25
$ws.Range("L32").Value = 9497.25
$ws.Range("N32").Value = -10149.25

# Row 101
$ws.Range("H101").Value = 1690
$ws.Range("J101").Value = 1700
$ws.Range("L101").Value = 5100
$ws.Range("N101").Value = -8344

# Row 138
$ws.Range("H138").Value = 2672.0793
$ws.Range("J138").Value = 2935.0205
$ws.Range("L138").Value = 8805.0615
$ws.Range("N138").Value = -19085.0615

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 11209.667
$ws.Range("I32").Value = 11209.667
$ws.Range("K32").Value = 11209.667
$ws.Range("M32").Value = -10922.667

# Row 61
$ws.Range("H61").Value = 2262.8215
$ws.Range("I61").Value = 1648.375
$ws.Range("K61").Value = 1648.375
$ws.Range("M61").Value = -1436.375

# Row 132
$ws.Range("H132").Value = 2618.875
$ws.Range("I132").Value = 1916.2142
$ws.Range("K132").Value = 5748.642599999999
$ws.Range("M132").Value = -3218.642599999999

# Row 136
$ws.Range("H136").Value = 2262.8215
$ws.Range("I136").Value = 1648.375
$ws.Range("K136").Value = 4945.125
$ws.Range("M136").Value = -2395.125

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 644.5
$ws.Range("I22").Value = 633.7
$ws.Range("K22").Value = 633.7
$ws.Range("M22").Value = -460.7

# Row 64
$ws.Range("H64").Value = 1000
$ws.Range("I64").Value = 1000
$ws.Range("J64").Value = 1000
$ws.Range("K64").Value = 1000
$ws.Range("L64").Value = 1000
$ws.Range("M64").Value = -775
$ws.Range("N64").Value = -1450

# Row 67
$ws.Range("H67").Value = 1000
$ws.Range("I67").Value = 1000
$ws.Range("J67").Value = 1000
$ws.Range("K67").Value = 1000
$ws.Range("L67").Value = 1000
$ws.Range("M67").Value = -220
$ws.Range("N67").Value = -2560

# Row 94
$ws.Range("H94").Value = 425.33334
$ws.Range("I94").Value = 425.33334
$ws.Range("K94").Value = 425.33334
$ws.Range("M94").Value = 25.66665999999998

# Row 107
$ws.Range("H107").Value = 2221.6667
$ws.Range("I107").Value = 2221.5
$ws.Range("J107").Value = 2222
$ws.Range("K107").Value = 2221.5
$ws.Range("L107").Value = 2222
$ws.Range("M107").Value = -301.5
$ws.Range("N107").Value = -6062

# Row 134
$ws.Range("H134").Value = 2120.087
$ws.Range("J134").Value = 2940
$ws.Range("L134").Value = 8820
$ws.Range("N134").Value = -13890

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 23998
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

# Row 34
$ws.Range("H34").Value = 23998
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

# Row 58
$ws.Range("H58").Value = 2659.3333
$ws.Range("I58").Value = 2265
$ws.Range("K58").Value = 2265
$ws.Range("M58").Value = -2062

# Row 68
$ws.Range("H68").Value = 49995
$ws.Range("J68").Value = 49995
$ws.Range("L68").Value = 49995
$ws.Range("N68").Value = -51493

# Row 71
$ws.Range("H71").Value = 49995
$ws.Range("J71").Value = 49995
$ws.Range("L71").Value = 149985
$ws.Range("N71").Value = -157473

# Row 132
$ws.Range("H132").Value = 2512.3
$ws.Range("I132").Value = 1939.7693
$ws.Range("J132").Value = 3575.5715
$ws.Range("K132").Value = 5819.3079
$ws.Range("L132").Value = 10726.7145
$ws.Range("M132").Value = -3289.3079
$ws.Range("N132").Value = -15786.7145

# Row 134
$ws.Range("H134").Value = 2790
$ws.Range("I134").Value = 2855.1333
$ws.Range("K134").Value = 8565.3999
$ws.Range("M134").Value = -6030.3999

# Row 136
$ws.Range("H136").Value = 2659.3333
$ws.Range("I136").Value = 2265
$ws.Range("K136").Value = 6795
$ws.Range("M136").Value = -4245

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 295.75
$ws.Range("I12").Value = 250.8
$ws.Range("J12").Value = 370.66666
$ws.Range("K12").Value = 752.4000000000001
$ws.Range("L12").Value = 1111.99998
$ws.Range("M12").Value = -579.4000000000001
$ws.Range("N12").Value = -1457.99998

# Row 34
$ws.Range("H34").Value = 3775
$ws.Range("I34").Value = 1900
$ws.Range("K34").Value = 5700
$ws.Range("M34").Value = -5616

# Row 41
$ws.Range("H41").Value = 3225
$ws.Range("I41").Value = 2966.6667
$ws.Range("J41").Value = 4000
$ws.Range("K41").Value = 8900.000100000001
$ws.Range("L41").Value = 12000
$ws.Range("M41").Value = -8562.000100000001
$ws.Range("N41").Value = -12676

# Row 42
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

# Row 43
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

# Row 76
$ws.Range("H76").Value = 5000
$ws.Range("I76").Value = 5000
$ws.Range("K76").Value = 15000
$ws.Range("M76").Value = -14617

# Row 79
$ws.Range("H79").Value = 5000
$ws.Range("I79").Value = 5000
$ws.Range("K79").Value = 15000
$ws.Range("M79").Value = -13674

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1083.9412
$ws.Range("I102").Value = 1058
$ws.Range("J102").Value = 1499
$ws.Range("K102").Value = 1058
$ws.Range("L102").Value = 1499
$ws.Range("M102").Value = 564
$ws.Range("N102").Value = -4743

# Row 122
$ws.Range("H122").Value = 1958.3334
$ws.Range("I122").Value = 1943.5
$ws.Range("K122").Value = 5830.5
$ws.Range("M122").Value = -3380.5

# Row 126
$ws.Range("H126").Value = 1230.5
$ws.Range("I126").Value = 972
$ws.Range("J126").Value = 1799.2
$ws.Range("K126").Value = 2916
$ws.Range("L126").Value = 5397.6
$ws.Range("M126").Value = -446
$ws.Range("N126").Value = -10337.6

# Row 132
$ws.Range("H132").Value = 2947.2856
$ws.Range("I132").Value = 1544.6666
$ws.Range("K132").Value = 4633.9998
$ws.Range("M132").Value = -2103.9998

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2166.3333
$ws.Range("J22").Value = 2166.3333
$ws.Range("L22").Value = 2166.3333
$ws.Range("N22").Value = -2756.3333

# Row 27
$ws.Range("H27").Value = 2166.3333
$ws.Range("J27").Value = 2166.3333
$ws.Range("L27").Value = 2166.3333
$ws.Range("N27").Value = -2380.3333

# Row 93
$ws.Range("H93").Value = 2886.75
$ws.Range("I93").Value = 3516.1667
$ws.Range("J93").Value = 998.5
$ws.Range("K93").Value = 3516.1667
$ws.Range("L93").Value = 998.5
$ws.Range("M93").Value = -2268.1667
$ws.Range("N93").Value = -3494.5

# Row 136
$ws.Range("H136").Value = 8002989
$ws.Range("I136").Value = 8002989
$ws.Range("K136").Value = 24008967
$ws.Range("M136").Value = -24006417

$ws = $wb.Worksheets.Item("WVR")
# Row 29
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

# Row 132
$ws.Range("H132").Value = 3050.6287
$ws.Range("I132").Value = 2632.7
$ws.Range("K132").Value = 7898.099999999999
$ws.Range("M132").Value = -5368.099999999999

# Row 135
$ws.Range("H135").Value = 149980
$ws.Range("J135").Value = 149980
$ws.Range("L135").Value = 149980
$ws.Range("N135").Value = -160120
